# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date", "Correspond Handoff Datetime"
# and "Correspond Handback DateTime" timestamps that get refreshed each
# time the handback status report is (re)generated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first row (G2)
$overview.Range("G2").Value = "2016-08-30 07:09:19"

# zh-cn sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
$zhcn.Range("H2").Value = "2016-08-30 07:09:14"
$zhcn.Range("K2").Value = "2016-08-30 07:09:37"

# de-de sheet: Correspond Handback DateTime (K2)
$dede.Range("K2").Value = "2016-08-30 07:09:44"
